$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-36) holds the "Förändrad" date serial value.
# All of these currently hold 45677 (2025-01-20) and should become
# 45678 (2025-01-21).
for ($row = 2; $row -le 36; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45677) {
        $cell.Value = 45678
    }
}
